$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) store numeric-looking strings as literal text
# in the source workbook. Force text format before assigning so Excel does not
# auto-convert them to numbers (which would also drop significant trailing zeros).
$textCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","E25","E27","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","E44","E45","D46","E46","E47","D49","E49","E50","E51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "307.69"
$ws.Range("E2").Value = "1.59%"
$ws.Range("D3").Value = "39.48"
$ws.Range("E3").Value = "10.73%"
$ws.Range("D4").Value = "5.090"
$ws.Range("E4").Value = "1.09%"
$ws.Range("D5").Value = "0.08153"
$ws.Range("E5").Value = "3.21%"
$ws.Range("D6").Value = "1.974"
$ws.Range("E6").Value = "6.60%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "7.903"
$ws.Range("E7").Value = "1.50%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9282"
$ws.Range("E8").Value = "1.00%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1408"
$ws.Range("E9").Value = "4.81%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1946"
$ws.Range("E10").Value = "2.35%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.09247"
$ws.Range("E11").Value = "1.69%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03497"
$ws.Range("E12").Value = "0.64%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09863"
$ws.Range("E13").Value = "0.29%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001409"
$ws.Range("E14").Value = "0.06%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005900"
$ws.Range("E15").Value = "-3.94%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.945"
$ws.Range("E16").Value = "6.22%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.172"
$ws.Range("E17").Value = "1.66%"
$ws.Range("D19").Value = "0.3452"
$ws.Range("E19").Value = "0.35%"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").Value = "-3.02%"
$ws.Range("D21").Value = "4.812"
$ws.Range("E21").Value = "-6.84%"
$ws.Range("D22").Value = "0.2617"
$ws.Range("E22").Value = "19.39%"
$ws.Range("D23").Value = "0.04469"
$ws.Range("E23").Value = "1.36%"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").Value = "0.64%"
$ws.Range("E25").Value = "-9.60%"
$ws.Range("E27").Value = "0.05%"
$ws.Range("D39").Value = "0.02112"
$ws.Range("E39").Value = "9.02%"
$ws.Range("D40").Value = "0.05159"
$ws.Range("E40").Value = "1.54%"
$ws.Range("D41").Value = "0.007463"
$ws.Range("E41").Value = "-2.03%"
$ws.Range("D42").Value = "0.01013"
$ws.Range("E42").Value = "-0.38%"
$ws.Range("D43").Value = "0.1367"
$ws.Range("E43").Value = "1.74%"
$ws.Range("E44").Value = "-0.88%"
$ws.Range("E45").Value = "-4.88%"
$ws.Range("D46").Value = "0.00006313"
$ws.Range("E46").Value = "2.11%"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D49").Value = "0.001601"
$ws.Range("E49").Value = "-3.56%"
$ws.Range("E50").Value = "0.04%"
$ws.Range("E51").Value = "0.04%"

Write-Host "Applied $($textCells.Count) text-format resets and 87 cell updates"
